$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Investments")

# Update max weight constraints (column C) for AGG (row 3) and IVV (row 4)
# to effectively remove the constraint (set to ~1), same as BIL in row 2
$ws.Range("C3").Value = 0.99999000000000005
$ws.Range("C4").Value = 0.99999000000000005

# Align A7's formatting (EEM) with A3/A5 (no fill) so the redundant,
# visually-identical cell style is dropped
$ws.Range("A3").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the active selection to C8
$ws.Range("C8").Select()
